$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New shared strings / DOI values used below will be created automatically
# when we assign string values to cells (COM handles sharedStrings table).

# --- Data rows 597-602: Male/Female nominal formants (study 10.1121/1.1913429)
$rows597 = @(
    @(597, "Male",   18, 30, "nominalF1", 286,  50),
    @(598, "Male",   18, 30, "nominalF2", 792,  50),
    @(599, "Male",   18, 30, "nominalF3", 2128, 50),
    @(600, "Female", 18, 30, "nominalF1", 270,  25),
    @(601, "Female", 18, 30, "nominalF2", 807,  25),
    @(602, "Female", 18, 30, "nominalF3", 2364, 25)
)

foreach ($r in $rows597) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = "10.1121/1.1913429"
}

# --- Data rows 603-608: NA group with AVERAGE() formulas (study 10.1016/j.jcomdis.2015.10.007)
$ws.Cells.Item(603, 1).Value = "NA"
$ws.Cells.Item(603, 2).Value = 5
$ws.Cells.Item(603, 3).Value = 7.3
$ws.Cells.Item(603, 4).Value = "nominalF1"
$ws.Cells.Item(603, 5).Formula = "=AVERAGE(360,1044)"
$ws.Cells.Item(603, 6).Value = 7
$ws.Cells.Item(603, 7).Value = "10.1016/j.jcomdis.2015.10.007"

$ws.Cells.Item(604, 1).Value = "NA"
$ws.Cells.Item(604, 2).Value = 5
$ws.Cells.Item(604, 3).Value = 7.3
$ws.Cells.Item(604, 4).Value = "nominalF2"
$ws.Cells.Item(604, 5).Formula = "=AVERAGE(1189,2514)"
$ws.Cells.Item(604, 6).Value = 7
$ws.Cells.Item(604, 7).Value = "10.1016/j.jcomdis.2015.10.007"

$ws.Cells.Item(605, 1).Value = "NA"
$ws.Cells.Item(605, 2).Value = 6.1
$ws.Cells.Item(605, 3).Value = 7.9
$ws.Cells.Item(605, 4).Value = "nominalF1"
$ws.Cells.Item(605, 5).Formula = "=AVERAGE(396,916)"
$ws.Cells.Item(605, 6).Value = 8
$ws.Cells.Item(605, 7).Value = "10.1016/j.jcomdis.2015.10.007"

$ws.Cells.Item(606, 1).Value = "NA"
$ws.Cells.Item(606, 2).Value = 6.1
$ws.Cells.Item(606, 3).Value = 7.9
$ws.Cells.Item(606, 4).Value = "nominalF2"
$ws.Cells.Item(606, 5).Formula = "=AVERAGE(1299,2055)"
$ws.Cells.Item(606, 6).Value = 8
$ws.Cells.Item(606, 7).Value = "10.1016/j.jcomdis.2015.10.007"

$ws.Cells.Item(607, 1).Value = "NA"
$ws.Cells.Item(607, 2).Value = 5
$ws.Cells.Item(607, 3).Value = 7
$ws.Cells.Item(607, 4).Value = "nominalF1"
$ws.Cells.Item(607, 5).Formula = "=AVERAGE(360,1027)"
$ws.Cells.Item(607, 6).Value = 90
$ws.Cells.Item(607, 7).Value = "10.1016/j.jcomdis.2015.10.007"

$ws.Cells.Item(608, 1).Value = "NA"
$ws.Cells.Item(608, 2).Value = 5
$ws.Cells.Item(608, 3).Value = 7
$ws.Cells.Item(608, 4).Value = "nominalF2"
$ws.Cells.Item(608, 5).Formula = "=AVERAGE(1181,2090)"
$ws.Cells.Item(608, 6).Value = 90
$ws.Cells.Item(608, 7).Value = "10.1016/j.jcomdis.2015.10.007"

# --- Data rows 609-610: Male minF/maxF (study 10.1016/s0095-4470(19)31416-0 )
$ws.Cells.Item(609, 1).Value = "Male"
$ws.Cells.Item(609, 2).Value = 17
$ws.Cells.Item(609, 3).Value = 26
$ws.Cells.Item(609, 4).Value = "minF"
$ws.Cells.Item(609, 5).Value = 62
$ws.Cells.Item(609, 6).Value = 157
$ws.Cells.Item(609, 7).Value = "10.1016/s0095-4470(19)31416-0 "

$ws.Cells.Item(610, 1).Value = "Male"
$ws.Cells.Item(610, 2).Value = 17
$ws.Cells.Item(610, 3).Value = 26
$ws.Cells.Item(610, 4).Value = "maxF"
$ws.Cells.Item(610, 5).Value = 392
$ws.Cells.Item(610, 6).Value = 157
$ws.Cells.Item(610, 7).Value = "10.1016/s0095-4470(19)31416-0 "

# --- Apply styles to match the rest of the table: E/F columns use the
# thousands-style (style used elsewhere for numeric results), G column
# uses the hyperlink style.
$srcEF = $ws.Range("E596:F596")
$srcG = $ws.Range("G596")

$srcEF.Copy() | Out-Null
$ws.Range("E597:F610").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$srcG.Copy() | Out-Null
$ws.Range("G597:G610").PasteSpecial(-4122) | Out-Null    # xlPasteFormats

$excel.CutCopyMode = 0

# --- Hyperlinks for the new rows
$ws.Hyperlinks.Add($ws.Range("G597"), "https://doi.org/10.1121/1.1913429", "", "", $ws.Range("G597")) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G598:G602"), "https://doi.org/10.1121/1.1913429", "", "", $ws.Range("G598:G602")) | Out-Null

# --- Fix the shared-formula ref for B584 (it should only cover B584, not B584:C584)
$ws.Range("B584").Formula = "=6+4/12"

# --- Update the view so it matches the saved state (scrolled further down)
$ws.Application.ActiveWindow.ScrollRow = 580
$ws.Range("B589").Select() | Out-Null
